$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "url" Name/Value row (currently row 5) down to row 8, shifting
# username/password/ppapplicationurl rows up by one.
$ws.Range("A5").Value = "username"
$ws.Range("B5").Value = "ppmaster"
$ws.Range("A6").Value = "password"
$ws.Range("B6").Value = "35Ramrod!"
$ws.Range("A7").Value = "ppapplicationurl"
$ws.Range("B7").Value = "https://wdts-gateway-env.wdts.local:792/login"
$ws.Range("A8").Value = "url"
$ws.Range("B8").Value = "https://tableIP:790/login/table-ui"

$ws.Range("E12").Select()
